$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
  $ws.Rows.Item(52).Copy()
  $ws.Rows.Item(55).Insert()
  Write-Output "OK"
} catch {
  Write-Output "ERROR: $_"
}
